$d = $word.ActiveDocument

# 1. "Prototyping Labs Supervisor" -> "Prototyping Lab Supervisor" (signature block)
$d.Content.Find.Execute("Labs Supervisor", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Lab Supervisor", 2) | Out-Null

# 2. "Prototyping Labs at GIX" -> "Prototyping Lab at GIX" (header)
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("Prototyping Labs at GIX", $false, $false, $false, $false, $false, `
                                    $true, 1, $false, "Prototyping Lab at GIX", 2) | Out-Null
        }
    }
}

# 3. "Ensure all other students are clear of immediate work area." ->
#    split into three runs, wrapping "immediate" with grammar proofing marks
$d.Content.Find.Execute("Ensure all other students are clear of immediate work area.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Ensure all other students are clear of immediate work area.", 2) | Out-Null
